# Update the build/version string throughout the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    foreach ($cell in $usedRange.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldVersion)) {
            $cell.Value2 = $val.Replace($oldVersion, $newVersion)
        }
    }
}
